# Insert 9 missing daily history rows (2019-11-18 .. 2019-11-28) into the
# SERBADK (id 5279) price history sheet, just before the existing
# 2019-11-29 row. This pushes the old rows 677..749 down to 686..758 and
# grows the used range from A1:I749 to A1:I758.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 9 blank rows starting at row 677 (old row 677, the
# 2019-11-29 record, becomes row 686; everything below shifts the same
# amount).
$ws.Rows("677:685").Insert()

# id / name are constant for every record in this sheet.
$stockId   = "5279"
$stockName = "SERBADK"

# timestamp (epoch seconds), date text, open, high, low, close, volume
$newRecords = @(
    @(1574035200, "2019-11-18", 2.773, 2.773, 2.754, 2.76,  4252948),
    @(1574121600, "2019-11-19", 2.779, 2.825, 2.773, 2.779, 12181044),
    @(1574208000, "2019-11-20", 2.818, 2.844, 2.805, 2.818, 5872497),
    @(1574294400, "2019-11-21", 2.831, 2.851, 2.805, 2.812, 3796798),
    @(1574380800, "2019-11-22", 2.818, 2.831, 2.779, 2.779, 3256348),
    @(1574640000, "2019-11-25", 2.799, 2.812, 2.773, 2.786, 1941149),
    @(1574726400, "2019-11-26", 2.799, 2.857, 2.779, 2.805, 9500095),
    @(1574812800, "2019-11-27", 2.818, 2.851, 2.812, 2.818, 4621048),
    @(1574899200, "2019-11-28", 2.831, 2.844, 2.812, 2.812, 2570699)
)

$r = 677
foreach ($rec in $newRecords) {
    $ws.Cells.Item($r, 1).Value = $rec[0]            # A: timestamp (number)
    $ws.Cells.Item($r, 2).Value = "'" + $rec[1]       # B: date text, e.g. 2019-11-18
    $ws.Cells.Item($r, 3).Value = "'" + $stockId      # C: id text "5279"
    $ws.Cells.Item($r, 4).Value = "'" + $stockName    # D: name text "SERBADK"
    $ws.Cells.Item($r, 5).Value = $rec[2]             # E: open
    $ws.Cells.Item($r, 6).Value = $rec[3]             # F: high
    $ws.Cells.Item($r, 7).Value = $rec[4]             # G: low
    $ws.Cells.Item($r, 8).Value = $rec[5]             # H: close
    $ws.Cells.Item($r, 9).Value = $rec[6]             # I: vol
    $r = $r + 1
}
